$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix typo "Totol" -> "Total" in the totals row label
$ws.Range("B29").Value = "Total"

# Fill in the missing SUM totals for the D/E/F/G columns (mirrors C29's
# existing SUM(C22:C28) pattern)
$ws.Range("D29").Formula = "=SUM(D22:D28)"
$ws.Range("E29").Formula = "=SUM(E22:E28)"
$ws.Range("F29").Formula = "=SUM(F22:F28)"
$ws.Range("G29").Formula = "=SUM(G22:G28)"

# Move the active selection to B29 (the totals label cell)
$ws.Range("B29").Select() | Out-Null
